# Fruta / hortaliza, semanal
# Insert the latest weekly price observation for
# "Feria Lagunitas de Puerto Montt - Haba" as a new row right after the
# existing row 38, pushing all the following (older) observations down by
# one row (row 39 -> 40, ..., row 96 -> 97). The sheet's used range grows
# from A1:R96 to A1:R97.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 39:96 down to 40:97, leaving a fresh (blank) row 39 that
# inherits formatting from the row above (same as Excel's native
# Rows.Insert behaviour).
$ws.Rows("39:39").Insert()

# Populate the newly inserted row 39 with the new weekly record.
$ws.Cells.Item(39, 1).Value = 4
$ws.Cells.Item(39, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(39, 3).Value = "Los Lagos"
$ws.Cells.Item(39, 4).Value = 44799
$ws.Cells.Item(39, 5).Value = 10
$ws.Cells.Item(39, 6).Value = 100112026
$ws.Cells.Item(39, 7).Value = "Haba"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 80
$ws.Cells.Item(39, 11).Value = 15000
$ws.Cells.Item(39, 12).Value = 15000
$ws.Cells.Item(39, 13).Value = 15000
$ws.Cells.Item(39, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(39, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(39, 16).Value = 600
$ws.Cells.Item(39, 17).Value = 25
$ws.Cells.Item(39, 18).Value = "Hortaliza"
